$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(2, 7).Value = 17.95625533333333
$ws.Cells.Item(2, 8).Value = 53.86876599999999
$ws.Cells.Item(2, 9).Value = 0.05960074617816258
$ws.Cells.Item(2, 10).Value = 0.05960074617816258
$ws.Cells.Item(2, 13).Value = 3.390429
$ws.Cells.Item(2, 14).Value = 10.171287
$ws.Cells.Item(2, 15).Value = 0.173121426386348
$ws.Cells.Item(2, 16).Value = 0.173121426386348
$ws.Cells.Item(2, 17).Value = 60.87940881353799
$ws.Cells.Item(2, 18).Value = 547.914679321842
$ws.Cells.Item(2, 19).Value = 0.01031816619205419
$ws.Cells.Item(2, 20).Value = 0.01031816619205419
$ws.Cells.Item(3, 7).Value = 17.95625533333333
$ws.Cells.Item(3, 8).Value = 53.86876599999999
$ws.Cells.Item(3, 9).Value = 0.05960074617816258
$ws.Cells.Item(3, 10).Value = 0.05960074617816258
$ws.Cells.Item(3, 15).Value = 0.5936336753560868
$ws.Cells.Item(3, 16).Value = 0.5936336753560868
$ws.Cells.Item(3, 17).Value = 208.7555998229473
$ws.Cells.Item(3, 18).Value = 1878.800398406526
$ws.Cells.Item(3, 19).Value = 0.0353810100077079
$ws.Cells.Item(3, 20).Value = 0.0353810100077079
$ws.Cells.Item(4, 7).Value = 17.95625533333333
$ws.Cells.Item(4, 8).Value = 53.86876599999999
$ws.Cells.Item(4, 9).Value = 0.05960074617816258
$ws.Cells.Item(4, 10).Value = 0.05960074617816258
$ws.Cells.Item(4, 13).Value = 4.546141666666667
$ws.Cells.Item(4, 14).Value = 13.638425
$ws.Cells.Item(4, 15).Value = 0.2321342018628743
$ws.Cells.Item(4, 16).Value = 0.2321342018628743
$ws.Cells.Item(4, 17).Value = 81.63168054817223
$ws.Cells.Item(4, 18).Value = 734.68512493355
$ws.Cells.Item(4, 19).Value = 0.01383537164449952
$ws.Cells.Item(4, 20).Value = 0.01383537164449953
$ws.Cells.Item(5, 7).Value = 17.95625533333333
$ws.Cells.Item(5, 8).Value = 53.86876599999999
$ws.Cells.Item(5, 9).Value = 0.05960074617816258
$ws.Cells.Item(5, 10).Value = 0.05960074617816258
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.021752
$ws.Cells.Item(5, 14).Value = 0.06525600000000001
$ws.Cells.Item(5, 15).Value = 0.001110696394691009
$ws.Cells.Item(5, 16).Value = 0.001110696394691009
$ws.Cells.Item(5, 17).Value = 0.3905844660106667
$ws.Cells.Item(5, 18).Value = 3.515260194096
$ws.Cells.Item(5, 19).Value = 0.00006619833390097912
$ws.Cells.Item(5, 20).Value = 0.00006619833390097912
$ws.Cells.Item(6, 8).Value = 789.271408
$ws.Cells.Item(6, 9).Value = 0.8732549183303921
$ws.Cells.Item(6, 10).Value = 0.8732549183303921
$ws.Cells.Item(6, 13).Value = 3.390429
$ws.Cells.Item(6, 14).Value = 10.171287
$ws.Cells.Item(6, 15).Value = 0.173121426386348
$ws.Cells.Item(6, 16).Value = 0.173121426386348
$ws.Cells.Item(6, 17).Value = 891.9895568513439
$ws.Cells.Item(6, 18).Value = 8027.906011662095
$ws.Cells.Item(6, 19).Value = 0.1511791370602513
$ws.Cells.Item(6, 20).Value = 0.1511791370602514
$ws.Cells.Item(7, 8).Value = 789.271408
$ws.Cells.Item(7, 9).Value = 0.8732549183303921
$ws.Cells.Item(7, 10).Value = 0.8732549183303921
$ws.Cells.Item(7, 15).Value = 0.5936336753560868
$ws.Cells.Item(7, 16).Value = 0.5936336753560868
$ws.Cells.Item(7, 18).Value = 27527.70382379429
$ws.Cells.Item(7, 19).Value = 0.5183935266912501
$ws.Cells.Item(7, 20).Value = 0.5183935266912501
$ws.Cells.Item(8, 8).Value = 789.271408
$ws.Cells.Item(8, 9).Value = 0.8732549183303921
$ws.Cells.Item(8, 10).Value = 0.8732549183303921
$ws.Cells.Item(8, 13).Value = 4.546141666666667
$ws.Cells.Item(8, 14).Value = 13.638425
$ws.Cells.Item(8, 15).Value = 0.2321342018628743
$ws.Cells.Item(8, 16).Value = 0.2321342018628743
$ws.Cells.Item(8, 17).Value = 1196.046544739156
$ws.Cells.Item(8, 18).Value = 10764.4189026524
$ws.Cells.Item(8, 19).Value = 0.202712333489455
$ws.Cells.Item(8, 20).Value = 0.2027123334894551
$ws.Cells.Item(9, 8).Value = 789.271408
$ws.Cells.Item(9, 9).Value = 0.8732549183303921
$ws.Cells.Item(9, 10).Value = 0.8732549183303921
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.021752
$ws.Cells.Item(9, 14).Value = 0.06525600000000001
$ws.Cells.Item(9, 15).Value = 0.001110696394691009
$ws.Cells.Item(9, 16).Value = 0.001110696394691009
$ws.Cells.Item(9, 17).Value = 5.722743888938667
$ws.Cells.Item(9, 18).Value = 51.504695000448
$ws.Cells.Item(9, 19).Value = 0.0009699210894357581
$ws.Cells.Item(9, 20).Value = 0.0009699210894357581
$ws.Cells.Item(10, 7).Value = 13.13303333333333
$ws.Cells.Item(10, 8).Value = 39.3991
$ws.Cells.Item(10, 9).Value = 0.04359141545488614
$ws.Cells.Item(10, 10).Value = 0.04359141545488615
$ws.Cells.Item(10, 13).Value = 3.390429
$ws.Cells.Item(10, 14).Value = 10.171287
$ws.Cells.Item(10, 15).Value = 0.173121426386348
$ws.Cells.Item(10, 16).Value = 0.173121426386348
$ws.Cells.Item(10, 17).Value = 44.52661707129999
$ws.Cells.Item(10, 18).Value = 400.7395536417
$ws.Cells.Item(10, 19).Value = 0.007546608021749784
$ws.Cells.Item(10, 20).Value = 0.007546608021749787
$ws.Cells.Item(11, 7).Value = 13.13303333333333
$ws.Cells.Item(11, 8).Value = 39.3991
$ws.Cells.Item(11, 9).Value = 0.04359141545488614
$ws.Cells.Item(11, 10).Value = 0.04359141545488615
$ws.Cells.Item(11, 15).Value = 0.5936336753560868
$ws.Cells.Item(11, 16).Value = 0.5936336753560868
$ws.Cells.Item(11, 17).Value = 152.6818481972333
$ws.Cells.Item(11, 18).Value = 1374.1366337751
$ws.Cells.Item(11, 19).Value = 0.02587733217045819
$ws.Cells.Item(11, 20).Value = 0.02587733217045819
$ws.Cells.Item(12, 7).Value = 13.13303333333333
$ws.Cells.Item(12, 8).Value = 39.3991
$ws.Cells.Item(12, 9).Value = 0.04359141545488614
$ws.Cells.Item(12, 10).Value = 0.04359141545488615
$ws.Cells.Item(12, 13).Value = 4.546141666666667
$ws.Cells.Item(12, 14).Value = 13.638425
$ws.Cells.Item(12, 15).Value = 0.2321342018628743
$ws.Cells.Item(12, 16).Value = 0.2321342018628743
$ws.Cells.Item(12, 17).Value = 59.70463004638889
$ws.Cells.Item(12, 18).Value = 537.3416704175
$ws.Cells.Item(12, 19).Value = 0.01011905843469296
$ws.Cells.Item(12, 20).Value = 0.01011905843469296
$ws.Cells.Item(13, 7).Value = 13.13303333333333
$ws.Cells.Item(13, 8).Value = 39.3991
$ws.Cells.Item(13, 9).Value = 0.04359141545488614
$ws.Cells.Item(13, 10).Value = 0.04359141545488615
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 0.6666666666666666
$ws.Cells.Item(13, 13).Value = 0.021752
$ws.Cells.Item(13, 14).Value = 0.06525600000000001
$ws.Cells.Item(13, 15).Value = 0.001110696394691009
$ws.Cells.Item(13, 16).Value = 0.001110696394691009
$ws.Cells.Item(13, 17).Value = 0.2856697410666667
$ws.Cells.Item(13, 18).Value = 2.5710276696
$ws.Cells.Item(13, 19).Value = 0.00004841682798521997
$ws.Cells.Item(13, 20).Value = 0.00004841682798521998
$ws.Cells.Item(14, 7).Value = 7.095921999999999
$ws.Cells.Item(14, 8).Value = 21.287766
$ws.Cells.Item(14, 9).Value = 0.02355292003655921
$ws.Cells.Item(14, 10).Value = 0.02355292003655921
$ws.Cells.Item(14, 13).Value = 3.390429
$ws.Cells.Item(14, 14).Value = 10.171287
$ws.Cells.Item(14, 15).Value = 0.173121426386348
$ws.Cells.Item(14, 16).Value = 0.173121426386348
$ws.Cells.Item(14, 17).Value = 24.058219730538
$ws.Cells.Item(14, 18).Value = 216.523977574842
$ws.Cells.Item(14, 19).Value = 0.004077515112292726
$ws.Cells.Item(14, 20).Value = 0.004077515112292727
$ws.Cells.Item(15, 7).Value = 7.095921999999999
$ws.Cells.Item(15, 8).Value = 21.287766
$ws.Cells.Item(15, 9).Value = 0.02355292003655921
$ws.Cells.Item(15, 10).Value = 0.02355292003655921
$ws.Cells.Item(15, 15).Value = 0.5936336753560868
$ws.Cells.Item(15, 16).Value = 0.5936336753560868
$ws.Cells.Item(15, 17).Value = 82.495677740614
$ws.Cells.Item(15, 18).Value = 742.4610996655259
$ws.Cells.Item(15, 19).Value = 0.01398180648667066
$ws.Cells.Item(15, 20).Value = 0.01398180648667066
$ws.Cells.Item(16, 7).Value = 7.095921999999999
$ws.Cells.Item(16, 8).Value = 21.287766
$ws.Cells.Item(16, 9).Value = 0.02355292003655921
$ws.Cells.Item(16, 10).Value = 0.02355292003655921
$ws.Cells.Item(16, 13).Value = 4.546141666666667
$ws.Cells.Item(16, 14).Value = 13.638425
$ws.Cells.Item(16, 15).Value = 0.2321342018628743
$ws.Cells.Item(16, 16).Value = 0.2321342018628743
$ws.Cells.Item(16, 17).Value = 32.25906666761666
$ws.Cells.Item(16, 18).Value = 290.33160000855
$ws.Cells.Item(16, 19).Value = 0.00546743829422677
$ws.Cells.Item(16, 20).Value = 0.005467438294226772
$ws.Cells.Item(17, 7).Value = 7.095921999999999
$ws.Cells.Item(17, 8).Value = 21.287766
$ws.Cells.Item(17, 9).Value = 0.02355292003655921
$ws.Cells.Item(17, 10).Value = 0.02355292003655921
$ws.Cells.Item(17, 11).Value = 2
$ws.Cells.Item(17, 12).Value = 0.6666666666666666
$ws.Cells.Item(17, 13).Value = 0.021752
$ws.Cells.Item(17, 14).Value = 0.06525600000000001
$ws.Cells.Item(17, 15).Value = 0.001110696394691009
$ws.Cells.Item(17, 16).Value = 0.001110696394691009
$ws.Cells.Item(17, 17).Value = 0.154350495344
$ws.Cells.Item(17, 18).Value = 1.389154458096
$ws.Cells.Item(17, 19).Value = 0.00002616014336905194
$ws.Cells.Item(17, 20).Value = 0.00002616014336905194
